$d = $word.ActiveDocument

# 1. Reword the "Click on the location..." bullet to the new text.
#    In the original document this sentence is split across two runs
#    ("C" + "lick on the location of an image and have it displayed on a
#    map."), but a plain text Find/Replace across the whole sentence
#    will match regardless of the run boundaries and leaves a single run.
$d.Content.Find.Execute(
    "Click on the location of an image and have it displayed on a map.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Select the name of the location of an image and open a map to display that location.",
    2) | Out-Null

# 2. Move the "_GoBack" bookmark from right after "Delete photos they have
#    already uploaded." to right before "Search for other users and view
#    their photos." (this mirrors where Word itself would redrop the
#    _GoBack bookmark after the edit above).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$searchRange = $d.Content
$searchRange.Find.Execute(
    "Search for other users and view their photos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$searchRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $searchRange) | Out-Null
